$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-11) were re-sorted. Capture current D,M,N,O,P,S values
# for every row first, then write them back per the new row order.
$cols = @("D", "M", "N", "O", "P", "S")
$original = @{}
foreach ($r in 2..11) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# Mapping: new row -> row whose original data now occupies it
$mapping = @{
    2  = 5
    3  = 3
    4  = 11
    5  = 2
    6  = 4
    7  = 7
    8  = 6
    9  = 10
    10 = 8
    11 = 9
}

foreach ($r in 2..11) {
    $src = $mapping[$r]
    $srcVals = $original[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcVals[$c]
    }
}
